# Update the cryptos price/volume table with refreshed figures.
# Rows 43/44 (Hedera <-> WhiteBITCoin) also swap places in this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.660.93"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.668.29"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.65"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.50"
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.124.06"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.573.00"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.42"
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.675.45"
$ws.Range("E16").Value = "  +6.69%  "
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.47"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.14"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.36"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.418"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0750"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.66"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.89"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.98"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.02"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.840"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.839"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.62"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "287.39"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.608"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0540"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.75"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.34"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.968.71"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.44"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.47"
$ws.Range("E51").Value = "  +0.13%  "
